$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the four new month columns (JUL, AUG, SEP, OCT) into both the
#    2021-22 block and the 2022-23 block, right before their existing
#    "Total" column. Excel's column insert shifts everything to the right
#    (values, merged cells, styles) automatically.
# ---------------------------------------------------------------------------

# Block 1 (2021-22): currently C,D,E = APR/MAY/JUN, F = Total.
# Insert 4 columns at F:I -> old F (Total) becomes J.
$ws.Range("F1:I1").EntireColumn.Insert()

# Block 2 (2022-23): after the previous insert it is K,L,M = APR/MAY/JUN,
# N = Total. Insert 4 columns at N:Q -> old N (Total) becomes R.
$ws.Range("N1:Q1").EntireColumn.Insert()

# The "Fig in crore" note in row 1 stays anchored at J1 in the source
# workbook (it is not part of the shifting month grid), so put it back at
# J1 and clear the R1 cell that the column insert shifted it into.
$ws.Range("R1").Value = $null
$ws.Range("J1").Value = "Fig in crore"

# ---------------------------------------------------------------------------
# 2. Fill in the header labels for the newly-inserted month columns.
# ---------------------------------------------------------------------------
$ws.Range("F3").Value = "JUL"
$ws.Range("G3").Value = "AUG"
$ws.Range("H3").Value = "SEP"
$ws.Range("I3").Value = "OCT"

$ws.Range("N3").Value = "JUL"
$ws.Range("O3").Value = "AUG"
$ws.Range("P3").Value = "SEP"
$ws.Range("Q3").Value = "OCT"

# ---------------------------------------------------------------------------
# 3. Update the monthly figures + totals for the four existing data rows
#    (PU1, PU10, PU11, PU15 -> rows 4-7) to the new values that include the
#    JUL-OCT months.
# ---------------------------------------------------------------------------

# Row 4 - PU1
$ws.Range("F4").Value = 216.22
$ws.Range("G4").Value = 215.2800000000001
$ws.Range("H4").Value = 214.7099999999998
$ws.Range("I4").Value = 216.9200000000001
$ws.Range("J4").Value = 1507.04
$ws.Range("N4").Value = 218.38
$ws.Range("O4").Value = 218.4200000000001
$ws.Range("P4").Value = 217.6899999999998
$ws.Range("Q4").Value = 211.6300000000001
$ws.Range("R4").Value = 1528.7

# Row 5 - PU10
$ws.Range("F5").Value = 17.92999999999999
$ws.Range("G5").Value = 20.22
$ws.Range("H5").Value = 21.35000000000001
$ws.Range("I5").Value = 18.70999999999999
$ws.Range("J5").Value = 131.63
$ws.Range("N5").Value = 21.55000000000001
$ws.Range("O5").Value = 22.29999999999998
$ws.Range("P5").Value = 22.36000000000001
$ws.Range("Q5").Value = 19.38
$ws.Range("R5").Value = 169.97

# Row 6 - PU11
$ws.Range("F6").Value = 0.1200000000000001
$ws.Range("G6").Value = 0.21
$ws.Range("H6").Value = 0.1899999999999999
$ws.Range("I6").Value = 0.04000000000000004
$ws.Range("J6").Value = 0.99
$ws.Range("N6").Value = 2.92
$ws.Range("O6").Value = 0.5800000000000001
$ws.Range("P6").Value = 2.77
$ws.Range("Q6").Value = 0.2900000000000009
$ws.Range("R6").Value = 9.880000000000001

# Row 7 - PU15
$ws.Range("F7").Value = 0.5700000000000001
$ws.Range("G7").Value = 0.7500000000000002
$ws.Range("H7").Value = 0.6099999999999999
$ws.Range("I7").Value = 0.73
$ws.Range("J7").Value = 4.04
$ws.Range("N7").Value = 0.9100000000000001
$ws.Range("O7").Value = 0.7500000000000004
$ws.Range("P7").Value = 0.6999999999999993
$ws.Range("Q7").Value = 0.6100000000000003
$ws.Range("R7").Value = 5.5

# ---------------------------------------------------------------------------
# 4. Insert two new rows above the existing PU32 row (currently row 8) to
#    make room for the new PU31 and STAFF rows, then populate all three plus
#    two brand-new rows below (PU27, PU28).
# ---------------------------------------------------------------------------
$ws.Range("A8:A9").EntireRow.Insert()

# Row 8 - PU31 (new)
$ws.Range("A8").Value = "PU31"
$ws.Range("B8").Value = 20.66
$ws.Range("C8").Value = 1.23
$ws.Range("D8").Value = 2.56
$ws.Range("E8").Value = 1.88
$ws.Range("F8").Value = 2.15
$ws.Range("G8").Value = 1.469999999999999
$ws.Range("H8").Value = 1.620000000000001
$ws.Range("I8").Value = 1.73
$ws.Range("J8").Value = 12.64
$ws.Range("K8").Value = 2.08
$ws.Range("L8").Value = 2.71
$ws.Range("M8").Value = 2.32
$ws.Range("N8").Value = 2.87
$ws.Range("O8").Value = 0.5
$ws.Range("P8").Value = 6.91
$ws.Range("Q8").Value = 2.989999999999998
$ws.Range("R8").Value = 20.38

# Row 9 - STAFF (new)
$ws.Range("A9").Value = "STAFF"
$ws.Range("B9").Value = 5075
$ws.Range("C9").Value = 407.59
$ws.Range("D9").Value = 378.6100000000001
$ws.Range("E9").Value = 403.4300000000001
$ws.Range("F9").Value = 455
$ws.Range("G9").Value = 421.0099999999998
$ws.Range("H9").Value = 424.7800000000002
$ws.Range("I9").Value = 540.29
$ws.Range("J9").Value = 3030.71
$ws.Range("K9").Value = 539.83
$ws.Range("L9").Value = 466.36
$ws.Range("M9").Value = 451.8199999999999
$ws.Range("N9").Value = 476.55
$ws.Range("O9").Value = 446.9099999999999
$ws.Range("P9").Value = 448.9200000000001
$ws.Range("Q9").Value = 566.02
$ws.Range("R9").Value = 3396.41

# Row 10 - PU32 (pre-existing row, shifted down; update with new month values)
$ws.Range("F10").Value = 30.8
$ws.Range("G10").Value = 31.2
$ws.Range("H10").Value = 29.72999999999999
$ws.Range("I10").Value = 38.14000000000001
$ws.Range("J10").Value = 209.56
$ws.Range("N10").Value = 29.19000000000001
$ws.Range("O10").Value = 30.88999999999999
$ws.Range("P10").Value = 33.78
$ws.Range("Q10").Value = 33.65000000000001
$ws.Range("R10").Value = 230.28

# Row 11 - PU27 (new)
$ws.Range("A11").Value = "PU27"
$ws.Range("B11").Value = 187.82
$ws.Range("C11").Value = 26.54
$ws.Range("D11").Value = 11.66
$ws.Range("E11").Value = 66.31999999999999
$ws.Range("F11").Value = -31.77
$ws.Range("G11").Value = 3.079999999999998
$ws.Range("H11").Value = 12.18000000000001
$ws.Range("I11").Value = 6.459999999999994
$ws.Range("J11").Value = 94.47
$ws.Range("K11").Value = 12.13
$ws.Range("L11").Value = 20.73
$ws.Range("M11").Value = 16.47
$ws.Range("N11").Value = 16.26000000000001
$ws.Range("O11").Value = 27.89
$ws.Range("P11").Value = 19.88
$ws.Range("Q11").Value = 10.33
$ws.Range("R11").Value = 123.69

# Row 12 - PU28 (new)
$ws.Range("A12").Value = "PU28"
$ws.Range("B12").Value = 113.64
$ws.Range("C12").Value = 6.03
$ws.Range("D12").Value = 4.29
$ws.Range("E12").Value = 7.539999999999999
$ws.Range("F12").Value = 6.949999999999999
$ws.Range("G12").Value = 10.85
$ws.Range("H12").Value = 10.38
$ws.Range("I12").Value = 6.82
$ws.Range("J12").Value = 52.86
$ws.Range("K12").Value = 23.68
$ws.Range("L12").Value = -2.359999999999999
$ws.Range("M12").Value = 9.27
$ws.Range("N12").Value = 9.180000000000003
$ws.Range("O12").Value = 5.989999999999995
$ws.Range("P12").Value = 13.57
$ws.Range("Q12").Value = 9.920000000000002
$ws.Range("R12").Value = 69.25
